$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment each date in F2:F7 by one day (shift forward by 1 day)
foreach ($row in 2..7) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 + 1
}
